$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Materials")

$row1Vals = @('Taxon_Local_ID','typeStatus','occurrenceID','catalogNumber','occurrenceDetails','occurrenceRemarks','recordNumber','recordedBy','individualID','individualCount','sex','lifeStage','reproductiveCondition','behavior','establishmentMeans','occurrenceStatus','preparations','disposition','otherCatalogNumbers','previousIdentifications','associatedMedia','associatedReferences','associatedOccurrences','associatedSequences','taxonID','scientificNameID','acceptedNameUsageID','parentNameUsageID','originalNameUsageID','nameAccordingToID','namePublishedInID','taxonConceptID','scientificName','acceptedNameUsage','parentNameUsage','originalNameUsage','nameAccordingTo','namePublishedIn','higherClassification','kingdom','phylum','class','order','suborder','infraorder','superfamily','family','genus','subgenus','specificEpithet','infraspecificEpithet','taxonRank','verbatimTaxonRank','scientificNameAuthorship','vernacularName','nomenclaturalCode','taxonomicStatus','nomenclaturalStatus','taxonRemarks','locationID','higherGeographyID','higherGeography','continent','waterBody','islandGroup','island','country','countryCode','stateProvince','county','municipality','locality','verbatimLocality','verbatimElevation','minimumElevationInMeters','maximumElevationInMeters','verbatimDepth','minimumDepthInMeters','maximumDepthInMeters','minimumDistanceAboveSurfaceInMeters','maximumDistanceAboveSurfaceInMeters','locationAccordingTo','locationRemarks','verbatimCoordinates','verbatimLatitude','verbatimLongitude','verbatimCoordinateSystem','verbatimSRS','decimalLatitude','decimalLongitude','geodeticDatum','coordinateUncertaintyInMeters','coordinatePrecision','pointRadiusSpatialFit','footprintWKT','footprintSRS','footprintSpatialFit','georeferencedBy','georeferenceProtocol','georeferenceSources','georeferenceVerificationStatus','georeferenceRemarks','identificationID','identifiedBy','dateIdentified','identificationReferences','identificationRemarks','identificationQualifier','geologicalContextID','earliestEonOrLowestEonothem','latestEonOrHighestEonothem','earliestEraOrLowestErathem','latestEraOrHighestErathem','earliestPeriodOrLowestSystem','latestPeriodOrHighestSystem','earliestEpochOrLowestSeries','latestEpochOrHighestSeries','earliestAgeOrLowestStage','latestAgeOrHighestStage','lowestBiostratigraphicZone','highestBiostratigraphicZone','lithostratigraphicTerms','group','formation','member','bed','eventID','samplingProtocol','samplingEffort','eventDate','eventTime','startDayOfYear','endDayOfYear','year','month','day','verbatimEventDate','habitat','fieldNumber','fieldNotes','eventRemarks','type','modified','language','rights','rightsHolder','accessRights','bibliographicCitation','institutionID','collectionID','datasetID','institutionCode','collectionCode','datasetName','ownerInstitutionCode','basisOfRecord','informationWithheld','dataGeneralizations','dynamicProperties','source')
$row2Vals = @('${iNaturalistTaxonId}','','${gbifOccurrenceID}','${catalogNumber}','','','','${recordedBy}','','${individualCount}','${sex}','','','','','PRESENT','','','','${taxonName}','','!references.ref','','','iNaturalist:${iNaturalistTaxonId}','','','','','','','','${summary.taxonName}','','','','','','','${kingdom}','${phylum}','${class}','${order}','${suborder}','${infraorder}','${superfamily}','${family}','${genus}','${subgenus}','','','','','${summary.Author}','','','','','','','','','','','','Galiano Island','Canada','CA','British Columbia','','','${locality}','','','','','${depth}','','','','','','','','','','','','${latitude}','${longitude}','','${coordinateUncertaintyInMeters}','','','','','','${recordedBy}','','','','','','','','','Confidence: ${confidence}','','','','','','','','','','','','','','','','','','','','','!references.samplingProtocol','','!Date:YYYY-MM-DD','!Date:HH:mm:ss','','','!Date:YYYY','!Date:M','!Date:D','','','${fieldNumber}','${fieldNotes}','${eventRemarks}','','','','','','','!references.ref','','','','!references.institutionCode','!references.collectionCode','!references.datasetName','','!references.basisOfRecord','','','','')

$n = $row1Vals.Length

$arr1 = New-Object 'object[,]' 1,$n
$arr2 = New-Object 'object[,]' 1,$n
for ($i = 0; $i -lt $n; $i++) {
    $arr1[0,$i] = $row1Vals[$i]
    $arr2[0,$i] = $row2Vals[$i]
}

$ws.Range("A1").Resize(1, $n).Value = $arr1
$ws.Range("A2").Resize(1, $n).Value = $arr2

Write-Output "done"
